$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AF2").Value = 15
$ws.Range("G2").Value = 2.45
$ws.Range("I2").Value = 3.4
$ws.Range("L2").Value = 1.8
$ws.Range("M2").Value = 1.91
$ws.Range("U2").Value = 9.5
$ws.Range("W2").Value = 23
# Row 4
$ws.Range("AE4").Value = 7.5
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.75
$ws.Range("V4").Value = 11
$ws.Range("Z4").Value = 7.5
# Row 5
$ws.Range("AC5").Value = 41
$ws.Range("AD5").Value = 126
$ws.Range("AE5").Value = 10
$ws.Range("Y5").Value = 29
# Row 6
$ws.Range("AD6").Value = 151
$ws.Range("AE6").Value = 9
$ws.Range("AF6").Value = 9.5
$ws.Range("G6").Value = 4.2
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 1.75
$ws.Range("N6").Value = 1.65
$ws.Range("O6").Value = 2.2
$ws.Range("R6").Value = 1.62
$ws.Range("S6").Value = 2.2
# Row 7
$ws.Range("AE7").Value = 12
$ws.Range("G7").Value = 2.1
$ws.Range("I7").Value = 3.25
$ws.Range("K7").Value = 13
$ws.Range("W7").Value = 19
$ws.Range("Z7").Value = 13
# Row 8
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 101
$ws.Range("AF8").Value = 17
$ws.Range("AG8").Value = 15
$ws.Range("AH8").Value = 41
$ws.Range("AI8").Value = 41
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 1.17
$ws.Range("K8").Value = 4.75
$ws.Range("L8").Value = 1.67
$ws.Range("M8").Value = 2.1
$ws.Range("N8").Value = 3.4
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 2.08
$ws.Range("R8").Value = 2.5
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 5
$ws.Range("U8").Value = 9
$ws.Range("W8").Value = 21
$ws.Range("X8").Value = 26
$ws.Range("Y8").Value = 51
$ws.Range("Z8").Value = 4.75
# Row 10
$ws.Range("AH10").Value = 51
$ws.Range("G10").Value = 1.91
$ws.Range("I10").Value = 4.5
$ws.Range("T10").Value = 5.5
$ws.Range("W10").Value = 15
# Row 11
$ws.Range("AB11").Value = 17
$ws.Range("AF11").Value = 19
$ws.Range("G11").Value = 1.91
$ws.Range("I11").Value = 3.75
$ws.Range("N11").Value = 2.15
$ws.Range("O11").Value = 1.67
$ws.Range("P11").Value = 1.5
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.73
$ws.Range("T11").Value = 6.5
$ws.Range("U11").Value = 8.5
$ws.Range("W11").Value = 17
$ws.Range("X11").Value = 17
$ws.Range("Y11").Value = 34
$ws.Range("Z11").Value = 8
# Row 12
$ws.Range("J12").Value = 1.04
$ws.Range("K12").Value = 13
# Row 13
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9
$ws.Range("N13").Value = 2.25
$ws.Range("O13").Value = 1.62
# Row 15
$ws.Range("AG15").Value = 9
$ws.Range("G15").Value = 3.2
$ws.Range("I15").Value = 2.2
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 1.91
$ws.Range("T15").Value = 10
$ws.Range("V15").Value = 12
$ws.Range("X15").Value = 26
# Row 17
$ws.Range("L17").Value = 1.21
$ws.Range("M17").Value = 3.95
$ws.Range("R17").Value = 1.6
$ws.Range("S17").Value = 2.19
$ws.Range("U17").Value = 8.5
# Row 18
$ws.Range("AF18").Value = 17
$ws.Range("AH18").Value = 34
$ws.Range("AI18").Value = 23
$ws.Range("G18").Value = 2.3
$ws.Range("I18").Value = 3.1
$ws.Range("K18").Value = 13
$ws.Range("W18").Value = 21
# Row 19
$ws.Range("AA19").Value = 8
$ws.Range("AF19").Value = 29
$ws.Range("AG19").Value = 19
$ws.Range("AI19").Value = 51
$ws.Range("G19").Value = 1.55
$ws.Range("H19").Value = 4.1
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 1.06
$ws.Range("K19").Value = 10
$ws.Range("L19").Value = 1.33
$ws.Range("M19").Value = 3.25
$ws.Range("N19").Value = 2.08
$ws.Range("O19").Value = 1.73
$ws.Range("T19").Value = 5.5
# Row 20
$ws.Range("AD20").Value = 351
$ws.Range("G20").Value = 2.7
$ws.Range("J20").Value = 1.07
$ws.Range("K20").Value = 9
$ws.Range("L20").Value = 1.4
$ws.Range("M20").Value = 2.75
$ws.Range("N20").Value = 2.25
$ws.Range("O20").Value = 1.62
# Row 21
$ws.Range("AD21").Value = 1000
$ws.Range("AF21").Value = 10
$ws.Range("AJ21").Value = 41
$ws.Range("J21").Value = 1.1
$ws.Range("K21").Value = 7
$ws.Range("N21").Value = 2.4
$ws.Range("O21").Value = 1.53
$ws.Range("T21").Value = 8
